$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 12675
$ws.Range("J2").Value = 300
$ws.Range("L2").Value = 300
$ws.Range("N2").Value = -526

$ws.Range("H4").Value = 59.833332
$ws.Range("I4").Value = 61.25
$ws.Range("K4").Value = 61.25
$ws.Range("M4").Value = 52.75

$ws.Range("H33").Value = 199.66667
$ws.Range("I33").Value = 198
$ws.Range("K33").Value = 198
$ws.Range("M33").Value = 31

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H70").Value = 4970.846
$ws.Range("J70").Value = 5714.6
$ws.Range("L70").Value = 17143.8
$ws.Range("N70").Value = -17683.8

$ws.Range("H73").Value = 4970.846
$ws.Range("J73").Value = 5714.6
$ws.Range("L73").Value = 17143.8
$ws.Range("N73").Value = -19015.8

$ws.Range("H97").Value = 561.5
$ws.Range("J97").Value = 561.5
$ws.Range("L97").Value = 1684.5
$ws.Range("N97").Value = -2676.5

$ws.Range("H112").Value = 2010.5714
$ws.Range("J112").Value = 1314.8
$ws.Range("L112").Value = 3944.4
$ws.Range("N112").Value = -6160.4

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H116").Value = 4196
$ws.Range("I116").Value = 4196
$ws.Range("K116").Value = 4196
$ws.Range("M116").Value = -754

$ws.Range("H118").Value = 846.3333
$ws.Range("I118").Value = 846.3333
$ws.Range("K118").Value = 2538.9999
$ws.Range("M118").Value = -881.9998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1707.3334
$ws.Range("I61").Value = 1707.3334
$ws.Range("K61").Value = 1707.3334
$ws.Range("M61").Value = -1495.3334

$ws.Range("H88").Value = 1331.4445
$ws.Range("I88").Value = 1362.3334
$ws.Range("K88").Value = 1362.3334
$ws.Range("M88").Value = -956.3334

$ws.Range("H91").Value = 1331.4445
$ws.Range("I91").Value = 1362.3334
$ws.Range("K91").Value = 1362.3334
$ws.Range("M91").Value = 41.66660000000002

$ws.Range("H132").Value = 2999.6667
$ws.Range("I132").Value = 2999.6667
$ws.Range("K132").Value = 8999.000100000001
$ws.Range("M132").Value = -6469.000100000001

$ws.Range("H136").Value = 1707.3334
$ws.Range("I136").Value = 1707.3334
$ws.Range("K136").Value = 5122.0002
$ws.Range("M136").Value = -2572.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4208.1665
$ws.Range("I86").Value = 1475
$ws.Range("J86").Value = 5574.75
$ws.Range("K86").Value = 1475
$ws.Range("L86").Value = 5574.75
$ws.Range("M86").Value = -352
$ws.Range("N86").Value = -7820.75

$ws.Range("H89").Value = 4208.1665
$ws.Range("I89").Value = 1475
$ws.Range("J89").Value = 5574.75
$ws.Range("K89").Value = 7375
$ws.Range("L89").Value = 27873.75
$ws.Range("M89").Value = -1759
$ws.Range("N89").Value = -39105.75

$ws.Range("H134").Value = 939.25
$ws.Range("I134").Value = 939.25
$ws.Range("K134").Value = 2817.75
$ws.Range("M134").Value = -282.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3835
$ws.Range("I86").Value = 3500
$ws.Range("K86").Value = 3500
$ws.Range("M86").Value = -2377

$ws.Range("H89").Value = 3835
$ws.Range("I89").Value = 3500
$ws.Range("K89").Value = 17500
$ws.Range("M89").Value = -11884

$ws.Range("H132").Value = 1661.7778
$ws.Range("J132").Value = 899.5
$ws.Range("L132").Value = 2698.5
$ws.Range("N132").Value = -7758.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 167.1875
$ws.Range("I12").Value = 27
$ws.Range("J12").Value = 230.90909
$ws.Range("K12").Value = 81
$ws.Range("L12").Value = 692.72727
$ws.Range("M12").Value = 92
$ws.Range("N12").Value = -1038.72727

$ws.Range("H14").Value = 786
$ws.Range("I14").Value = 786
$ws.Range("K14").Value = 2358
$ws.Range("M14").Value = -2185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13014.4375
$ws.Range("J43").Value = 19629.1
$ws.Range("L43").Value = 19629.1
$ws.Range("N43").Value = -19931.1

$ws.Range("H113").Value = 3285.4
$ws.Range("I113").Value = 1540.1666
$ws.Range("K113").Value = 1540.1666
$ws.Range("M113").Value = 629.8334

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 73192.07000000001
$ws.Range("I132").Value = 92508.09
$ws.Range("K132").Value = 277524.27
$ws.Range("M132").Value = -274994.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1000
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -830
$ws.Range("N16").Value = -2340

$ws.Range("H68").Value = 10000
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498

$ws.Range("H71").Value = 10000
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488

$ws.Range("H82").Value = 4799.8887
$ws.Range("I82").Value = 700
$ws.Range("K82").Value = 700
$ws.Range("M82").Value = -339

$ws.Range("H85").Value = 4799.8887
$ws.Range("I85").Value = 700
$ws.Range("K85").Value = 700
$ws.Range("M85").Value = 548

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1020.0769
$ws.Range("I113").Value = 964.6667
$ws.Range("J113").Value = 1144.75
$ws.Range("K113").Value = 2894.0001
$ws.Range("L113").Value = 3434.25
$ws.Range("M113").Value = -724.0001000000002
$ws.Range("N113").Value = -7774.25

$ws.Range("H132").Value = 1097.8667
$ws.Range("I132").Value = 1097.8667
$ws.Range("K132").Value = 3293.6001
$ws.Range("M132").Value = -763.6001000000001

$ws.Range("H136").Value = 2321.8276
$ws.Range("I136").Value = 1974.08
$ws.Range("J136").Value = 4495.25
$ws.Range("K136").Value = 5922.24
$ws.Range("L136").Value = 13485.75
$ws.Range("M136").Value = -3372.24
$ws.Range("N136").Value = -18585.75
